$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 482, shifting existing rows (482-532) down to (483-533)
$ws.Rows.Item(482).Insert()

# Populate the new row 482 with its data
$ws.Cells.Item(482, 1).Value = 7
$ws.Cells.Item(482, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(482, 3).Value = "Ñuble"
$ws.Cells.Item(482, 4).Value = 45132
$ws.Cells.Item(482, 5).Value = 16
$ws.Cells.Item(482, 6).Value = 100112008
$ws.Cells.Item(482, 7).Value = "Coliflor"
$ws.Cells.Item(482, 8).Value = "Sin especificar"
$ws.Cells.Item(482, 9).Value = "Primera"
$ws.Cells.Item(482, 10).Value = 250
$ws.Cells.Item(482, 11).Value = 1000
$ws.Cells.Item(482, 12).Value = 1000
$ws.Cells.Item(482, 13).Value = 1000
$ws.Cells.Item(482, 14).Value = "$/unidad"
$ws.Cells.Item(482, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(482, 16).Value = 1000
$ws.Cells.Item(482, 17).Value = 1
$ws.Cells.Item(482, 18).Value = "Hortaliza"
